$d = $word.ActiveDocument

# Locate the "Author" styled paragraph that holds "Edison Achalma" (the
# author name line right under the title) so we can add the affiliation
# paragraph right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Author" -and $p.Range.Text.Trim() -eq "Edison Achalma") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a brand-new paragraph right after the end of the target
    # paragraph (its Range.End already sits just past the paragraph mark,
    # i.e. at the start of the following paragraph), so the new paragraph
    # is created between "Edison Achalma" and whatever came next.
    $insertionPoint = $d.Range($target.Range.End, $target.Range.End)
    $insertionPoint.InsertParagraphAfter()

    # The freshly created paragraph inherited the style/formatting of the
    # paragraph that used to follow "Edison Achalma" -- force it back to
    # the "Author" style and fill in the affiliation text.
    $newPara = $target.Next()
    $newPara.Style = "Author"
    $newPara.Range.InsertAfter("Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga")
}
